$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; this shifts the existing rows 10-20
# down to rows 11-21 (and copies down the date style on column D).
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with the new price-report entry.
$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "Vega Modelo de Temuco"
$ws.Range("C10").Value = "La Araucanía"
$ws.Range("D10").Value = 45219
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100104
$ws.Range("H10").Value = "Frutos de pepita"
$ws.Range("I10").Value = 100104004
$ws.Range("J10").Value = "Níspero"
$ws.Range("K10").Value = "Californiana(o)"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 20
$ws.Range("N10").Value = 35000
$ws.Range("O10").Value = 35000
$ws.Range("P10").Value = 35000
$ws.Range("Q10").Value = "$/bandeja 10 kilos"
$ws.Range("R10").Value = "Provincia de Los Andes"
$ws.Range("S10").Value = 3500
$ws.Range("T10").Value = 10
